$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Fix styles on cells that lose their (no-op) fill flag ----
# G3, G11, G12 should match the plain style used by G2 (no border/fill)
$ws.Range("G2").Copy() | Out-Null
$ws.Range("G3").PasteSpecial(-4122) | Out-Null
$ws.Range("G11").PasteSpecial(-4122) | Out-Null
$ws.Range("G12").PasteSpecial(-4122) | Out-Null

# B12, F12 should match the style used by B2/F2 (left+right border)
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B12").PasteSpecial(-4122) | Out-Null
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F12").PasteSpecial(-4122) | Out-Null

# E12 should match the style used by E2 (right border)
$ws.Range("E2").Copy() | Out-Null
$ws.Range("E12").PasteSpecial(-4122) | Out-Null

# ---- Build new row 13 by cloning formats of row 12 (now clean) ----
$ws.Range("A12:G12").Copy() | Out-Null
$ws.Range("A13:G13").PasteSpecial(-4122) | Out-Null
# C13 needs the numeric-format style used by C7/C9/C11
$ws.Range("C7").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---- Write the refreshed variable table (values) ----
# Row 2
$ws.Range("A2").Value2 = 'uncert_tree_vulnerability_mean'
$ws.Range("B2").Value2 = 'tnorm_0_1'
$ws.Range("C2").Value2 = 0.1
$ws.Range("D2").Value2 = 'NA'
$ws.Range("E2").Value2 = 0.3
$ws.Range("F2").Value2 = 'relative'
$ws.Range("G2").Value2 = 'Unknown factors which make the trees vulnerable'

# Row 3
$ws.Range("A3").Value2 = 'uncert_tree_vulnerability_var'
$ws.Range("B3").Value2 = 'posnorm'
$ws.Range("C3").Value2 = 10
$ws.Range("D3").Value2 = 'NA'
$ws.Range("E3").Value2 = 20
$ws.Range("F3").Value2 = 'percent'
$ws.Range("G3").Value2 = 'Coefficient of variation of unknown factors for vulnerability'

# Row 4
$ws.Range("A4").Value2 = 'uncert_tree_parameter_age_1'
$ws.Range("B4").Value2 = 'posnorm'
$ws.Range("C4").Value2 = 20
$ws.Range("D4").Value2 = 'NA'
$ws.Range("E4").Value2 = 40
$ws.Range("F4").Value2 = '-'
$ws.Range("G4").Value2 = '"Best ages", turning point of curve on y-axis'

# Row 5
$ws.Range("A5").Value2 = 'uncert_tree_parameter_age_2'
$ws.Range("B5").Value2 = 'posnorm'
$ws.Range("C5").Value2 = 5
$ws.Range("D5").Value2 = 'NA'
$ws.Range("E5").Value2 = 7
$ws.Range("F5").Value2 = '-'
$ws.Range("G5").Value2 = 'Range of vulnerability, "curve compression"'

# Row 6
$ws.Range("A6").Value2 = 'uncert_influence_quali'
$ws.Range("B6").Value2 = 'tnorm_0_1'
$ws.Range("C6").Value2 = 0.1
$ws.Range("D6").Value2 = 'NA'
$ws.Range("E6").Value2 = 0.9
$ws.Range("F6").Value2 = 'percent per Eur'
$ws.Range("G6").Value2 = 'Unknown factor that shows the influence of fruit quality'

# Row 7
$ws.Range("A7").Value2 = 'uncert_good_direct_market_maximum_sells_kg'
$ws.Range("B7").Value2 = 'posnorm'
$ws.Range("C7").Value2 = 1500
$ws.Range("D7").Value2 = 'NA'
$ws.Range("E7").Value2 = 5000
$ws.Range("F7").Value2 = 'kg'
$ws.Range("G7").Value2 = 'How much kg wlanut the farmer could sell directly if have build a good supply chain'

# Row 8
$ws.Range("A8").Value2 = 'tree_labor_establishment_var_h_per_tree'
$ws.Range("B8").Value2 = 'posnorm'
$ws.Range("C8").Value2 = 20
$ws.Range("D8").Value2 = 'NA'
$ws.Range("E8").Value2 = 70
$ws.Range("F8").Value2 = 'percent'
$ws.Range("G8").Value2 = 'Coefficient of variation of minimum quality for direct marketing'

# Row 9
$ws.Range("A9").Value2 = 'uncert_poor_direct_market_maximum_sells_kg'
$ws.Range("B9").Value2 = 'posnorm'
$ws.Range("C9").Value2 = 1000
$ws.Range("D9").Value2 = 'NA'
$ws.Range("E9").Value2 = 1500
$ws.Range("F9").Value2 = 'kg'
$ws.Range("G9").Value2 = 'How much kg wlanut the farmer could sell directly if have build a poor supply chain'

# Row 10
$ws.Range("A10").Value2 = 'uncert_invest_until_good_market_mean_h'
$ws.Range("B10").Value2 = 'posnorm'
$ws.Range("C10").Value2 = 10
$ws.Range("D10").Value2 = 'NA'
$ws.Range("E10").Value2 = 20
$ws.Range("F10").Value2 = 'h'
$ws.Range("G10").Value2 = 'Unknown threshold: how much labor needs to be inested until enough network is build to enhance capapbilities of direct marketing'

# Row 11
$ws.Range("A11").Value2 = 'uncert_invest_until_good_market_var'
$ws.Range("B11").Value2 = 'posnorm'
$ws.Range("C11").Value2 = 10
$ws.Range("D11").Value2 = 'NA'
$ws.Range("E11").Value2 = 20
$ws.Range("F11").Value2 = 'percent'
$ws.Range("G11").Value2 = 'Coefficient of variation of unknown threshold of labor for a good local market'

# Row 12
$ws.Range("A12").Value2 = 'uncert_minimum_quali_for_direct_percent'
$ws.Range("B12").Value2 = 'posnorm'
$ws.Range("C12").Value2 = 0.7
$ws.Range("D12").Value2 = 'NA'
$ws.Range("E12").Value2 = 0.8
$ws.Range("F12").Value2 = 'relative'
$ws.Range("G12").Value2 = 'Minimum fruit quality [relative] so they are marketable directly'

# Row 13
$ws.Range("A13").Value2 = 'uncert_hay_good_market_capacity'
$ws.Range("B13").Value2 = 'posnorm'
$ws.Range("C13").Value2 = 3
$ws.Range("D13").Value2 = 'NA'
$ws.Range("E13").Value2 = 10
$ws.Range("F13").Value2 = 't hay per ha'
$ws.Range("G13").Value2 = 'Hay yield that indicates if overall, a lot of hay is available at market - low prices! (conceptual, uncertain - thus a range between possible maximum yields)'

# ---- Misc sheet metadata ----
$ws.Columns.Item(1).ColumnWidth = 48
$ws.Range("E3").Select() | Out-Null
